$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Values mirror the site's
# scraped price/volume formatting (plain decimal / percent strings),
# which must be written back as literal TEXT, not auto-converted numbers.
$updates = [ordered]@{
    "D2" = "273.60"
    "E2" = "1.30%"
    "D3" = "26.83"
    "E3" = "0.47%"
    "D4" = "4.914"
    "E4" = "4.41%"
    "D5" = "0.06331"
    "E5" = "3.88%"
    "D6" = "6.946"
    "E6" = "3.02%"
    "D7" = "3.353"
    "E7" = "5.86%"
    "D8" = "1.367"
    "E8" = "53.08%"
    "D9" = "0.8881"
    "E9" = "3.47%"
    "D10" = "0.1472"
    "E10" = "3.65%"
    "D11" = "0.05121"
    "E11" = "2.12%"
    "D12" = "0.07347"
    "E12" = "3.35%"
    "D13" = "0.03164"
    "E13" = "-0.53%"
    "D14" = "0.09053"
    "E14" = "0.28%"
    "D15" = "0.001559"
    "E15" = "1.87%"
    "D16" = "0.0006343"
    "E16" = "4.77%"
    "D17" = "0.006028"
    "E17" = "-0.87%"
    "D18" = "3.479"
    "E18" = "0.48%"
    "E19" = "1.67%"
    "E20" = "2.20%"
    "D21" = "0.1333"
    "E21" = "2.52%"
    "D22" = "3.919"
    "E22" = "2.25%"
    "D23" = "0.04345"
    "E23" = "2.45%"
    "D24" = "0.001181"
    "E24" = "-0.35%"
    "D25" = "0.003650"
    "E25" = "-12.03%"
    "D26" = "0.0001205"
    "E26" = "0.46%"
    "D27" = "0.0001941"
    "E27" = "15.51%"
    "D40" = "0.04029"
    "E40" = "1.90%"
    "D41" = "0.006620"
    "E41" = "57.94%"
    "D42" = "0.1166"
    "E42" = "4.24%"
    "D43" = "0.002371"
    "E43" = "17.96%"
    "D44" = "0.01260"
    "E44" = "-0.67%"
    "D45" = "0.00005252"
    "E45" = "2.50%"
    "D47" = "0.02125"
    "E47" = "-13.19%"
    "E48" = "-0.05%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage: a plain .Value = "274.38"/"1.57%" would be
    # auto-coerced to a Double (and the percent string to a percentage-
    # formatted fraction), same as typing it into a General cell in Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop the Text number format again so the cell keeps its original
    # (unstyled) appearance -- only the content changes, per the diff.
    $cell.ClearFormats()
}
